# Re-shuffle "Vennegruppe" group assignments on sheet 1, and append the
# newly-recorded group history entries (columns N:V) on sheet 2 ("Arkiv").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Vennegruppe": rows 2-6, columns A-E hold the current group
# membership. The names are the same set, just re-arranged into new
# groups, so we overwrite the cell values in place.
# ---------------------------------------------------------------------
$wsGroups = $wb.Worksheets.Item("Vennegruppe")

$wsGroups.Range("A2").Value = "Ria"
$wsGroups.Range("B2").Value = "Klara"
$wsGroups.Range("C2").Value = "Abel"
$wsGroups.Range("D2").Value = "Ulu"
$wsGroups.Range("E2").Value = "Iris"

$wsGroups.Range("A3").Value = "Vance"
$wsGroups.Range("B3").Value = "Dora"
$wsGroups.Range("C3").Value = "Lars"
$wsGroups.Range("D3").Value = "Tanya"
$wsGroups.Range("E3").Value = "Fransy"

$wsGroups.Range("A4").Value = "Penny"
$wsGroups.Range("B4").Value = "Maya"
$wsGroups.Range("C4").Value = "John"
$wsGroups.Range("D4").Value = "Benni"
$wsGroups.Range("E4").Value = "Stan"

$wsGroups.Range("A5").Value = "Wyatt"
$wsGroups.Range("B5").Value = "Q"
$wsGroups.Range("C5").Value = "Gert"
$wsGroups.Range("D5").Value = "Ede"
$wsGroups.Range("E5").Value = "Chris"

$wsGroups.Range("A6").Value = "Nick"
$wsGroups.Range("B6").Value = "Hermine"
$wsGroups.Range("C6").Value = "Olav"

# ---------------------------------------------------------------------
# Sheet "Arkiv": for each child (row 6-28), append the names of the new
# group members onto the end of their "Har vaert i gruppen med" history
# row, continuing right after whatever the last filled column already
# was.
# ---------------------------------------------------------------------
$wsArkiv = $wb.Worksheets.Item("Arkiv")

function Set-History($row, $values) {
    $lastCol = $wsArkiv.Cells.Item($row, 256).End(-4159).Column
    for ($i = 0; $i -lt $values.Count; $i++) {
        $wsArkiv.Cells.Item($row, $lastCol + 1 + $i).Value = $values[$i]
    }
}

Set-History 6  @("John","Lars","Stan","Gert","Lars","Q","Nick","Vance")
Set-History 7  @("Q","Olav","Nick","Ria","Stan","Klara","Penny")
Set-History 8  @("Klara","Tanya","Vance","John","Wyatt","Olav")
Set-History 9  @("Hermine","Ulu","Ede","Maya","Gert","Ede","Ulu","Hermine")
Set-History 10 @("Hermine","Ulu","Dora","Maya","Dora","Gert","Ulu","Hermine")
Set-History 11 @("Penny","Wyatt","Ria","Iris","Maya","Iris","Tanya")
Set-History 12 @("John","Abel","Lars","Stan","Dora","Ede","Ulu","Hermine")
Set-History 13 @("Ulu","Ede","Dora","Maya","Dora","Gert","Ede","Ulu")
Set-History 14 @("Penny","Wyatt","Ria","Fransy","Maya","Fransy","Tanya")
Set-History 15 @("Abel","Lars","Stan","Gert","Chris","Wyatt","Olav")
Set-History 16 @("Tanya","Chris","Vance","Ria","Benni","Stan","Penny")
Set-History 17 @("John","Abel","Stan","Gert","Abel","Q","Nick","Vance")
Set-History 18 @("Hermine","Ulu","Ede","Dora","Iris","Fransy","Tanya")
Set-History 19 @("Benni","Q","Olav","Lars","Abel","Q","Vance")
Set-History 20 @("Benni","Q","Nick","Chris","John","Wyatt")
Set-History 21 @("Wyatt","Ria","Iris","Fransy","Ria","Benni","Stan","Klara")
Set-History 22 @("Benni","Olav","Nick","Lars","Abel","Nick","Vance")
Set-History 23 @("Penny","Wyatt","Iris","Fransy","Benni","Stan","Klara","Penny")
Set-History 24 @("John","Abel","Lars","Gert","Ria","Benni","Klara","Penny")
Set-History 25 @("Klara","Chris","Vance","Maya","Iris","Fransy")
Set-History 26 @("Hermine","Ede","Dora","Maya","Dora","Gert","Ede","Hermine")
Set-History 27 @("Klara","Tanya","Chris","Lars","Abel","Q","Nick")
Set-History 28 @("Penny","Ria","Iris","Fransy","Chris","John","Olav")
